$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 10 (lastLogin), shifting rows 10-60 down to 11-61.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the "status" / "TEXT" pair (User Table section).
$ws.Range("A10").Value = "status"
$ws.Range("B10").Value = "TEXT"

# Move the active selection to B11 (the cell below the newly inserted row), matching the saved view state.
$ws.Range("B11").Select()
